$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for c1a15a61-... row (row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 19:08:20"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 19:08:15"
$wsZhCn.Range("K2").Value = "2016-08-25 19:08:33"

# de-de sheet: Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-25 19:08:40"
